$d = $word.ActiveDocument

# The last four paragraphs of the document (under "IMPORTANT") are being
# replaced with new content, and a new trailing empty paragraph is added.
#
# Old paragraph 30: "- Le probleme avec mes graphs c-est la saisonabilite, ..."
# Old paragraph 31: "C-est pour ca que tu as 2 piques dans ta simulation, ..."
# Old paragraph 32: (empty)
# Old paragraph 33: "Variable categorielle, comment evaluer les distances ..."
#
# New paragraph 30: "Prends une première période A avec un nombre fixe ..."
# New paragraph 31: "Puis compare labels période A  avec chacune pour ARI."
# New paragraph 32: (empty)
# New paragraph 33 (style BodyText): "tituscodes@gmail.com"
# New paragraph 34 (new, style Normal): (empty)

$count = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($count - 3)
$p2 = $d.Paragraphs.Item($count - 2)
$p4 = $d.Paragraphs.Item($count)

# Sanity-check we are pointed at the expected paragraphs before rewriting.
if ($p1.Range.Text -notmatch "Le probleme avec mes graphs") {
    throw "Unexpected content in target paragraph 1"
}
if ($p2.Range.Text -notmatch "2 piques dans ta simulation") {
    throw "Unexpected content in target paragraph 2"
}
if ($p4.Range.Text -notmatch "Variable categorielle") {
    throw "Unexpected content in target paragraph 4"
}

$p1.Range.Text = "Prends une première période A avec un nombre fixe de clients sans commande dans le futur. Puis itère, prend les clients période A + 1 week avec la meme méthode SQL et entraine a chaque fois."
$p2.Range.Text = "Puis compare labels période A  avec chacune pour ARI."

# third paragraph (empty) stays empty / untouched

$p4.Range.Text = "tituscodes@gmail.com"

# Add a fresh, empty, Normal-style paragraph after the last one while $p4
# is still "Normal" styled, so the new paragraph inherits/keeps an explicit
# Normal pStyle once touched.
$p4.Range.InsertParagraphAfter()
$newCount = $d.Paragraphs.Count
$p5 = $d.Paragraphs.Item($newCount)
$touch = $p5.Range.ParagraphStyle

# Now switch the email paragraph to the BodyText style.
$p4.Style = "Body Text"
